$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the report date (A2)
$ws.Range("A2").Value = 45639

# Update "Percentagem Execução Reportada Anteriormente" (F) and
# "Percentagem Execução Actual" (G) for tasks T3.1 - T3.2.3
$ws.Range("F21").Value = 0.7
$ws.Range("G21").Value = 0.8

$ws.Range("F22").Value = 0.8
$ws.Range("G22").Value = 0.9

$ws.Range("F23").Value = 0.65
$ws.Range("G23").Value = 0.8

$ws.Range("F24").Value = 0.6
$ws.Range("G24").Value = 0.8

$ws.Range("F25").Value = 0.6
$ws.Range("G25").Value = 0.8

$ws.Range("F26").Value = 0.35
$ws.Range("G26").Value = 0.65

$ws.Range("F27").Value = 0.7
$ws.Range("G27").Value = 0.9

$ws.Range("F28").Value = 0.3
$ws.Range("G28").Value = 0.5

$ws.Range("F29").Value = 0.3
$ws.Range("G29").Value = 0.5

# The "Reported previously" column picks up the same number format /
# border as the "Actual" column once both are filled in (xlPasteFormats).
$ws.Range("G21").Copy()
$ws.Range("F21").PasteSpecial(-4122)
$ws.Range("G26").Copy()
$ws.Range("F26").PasteSpecial(-4122)
$ws.Range("G27").Copy()
$ws.Range("F27").PasteSpecial(-4122)
$ws.Range("G28").Copy()
$ws.Range("F28").PasteSpecial(-4122)
$ws.Range("G29").Copy()
$ws.Range("F29").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Reset the view: scroll back to the top-left (remove the stale
# topLeftCell) and move the active selection to A2.
[void]$ws.Range("A2").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
